$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Forces a literal text value even when the text looks like a date
    # (avoids Excel auto-converting "1401-04-30"-style strings to date serials)
    $range.Formula = '=TEXT("' + $text + '","@")'
    $range.Copy()
    $range.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# --- Row 8: period headers (shift left one quarter, append Q4 1401/12 in M8) ---
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# --- Row 9: publish dates for each quarter column (kept as text) ---
$ws.Range("D9").Value = "1400-10-30 (2)"
$ws.Range("E9").Value = "1401-01-31 (8)"
$ws.Range("F9").Value = "1401-04-30 (2)"
$ws.Range("G9").Value = "1401-08-02 (4)"
$ws.Range("H9").Value = "1401-12-29 (3)"
$ws.Range("I9").Value = "1402-01-30 (9)"
Set-TextValue $ws.Range("J9") "1401-04-30"
$ws.Range("K9").Value = "1401-08-02 (2)"
$ws.Range("L9").Value = "1401-12-29 (2)"
$ws.Range("M9").Value = "1402-01-30 (2)"

# --- Numeric data rows 11-27: shift each quarter column left, append new quarter (M) values ---
# Row 11
$ws.Range("D11").Value = 1127846
$ws.Range("E11").Value = 745527
$ws.Range("F11").Value = 1882542
$ws.Range("G11").Value = 2769054
$ws.Range("H11").Value = 1798057
$ws.Range("I11").Value = 1056678
$ws.Range("J11").Value = 2922768
$ws.Range("K11").Value = 2650496
$ws.Range("L11").Value = 2371051
$ws.Range("M11").Value = 2271595
# Row 12
$ws.Range("D12").Value = -540034
$ws.Range("E12").Value = -439125
$ws.Range("F12").Value = -770656
$ws.Range("G12").Value = -887173
$ws.Range("H12").Value = -903526
$ws.Range("I12").Value = -694021
$ws.Range("J12").Value = -1190313
$ws.Range("K12").Value = -1055551
$ws.Range("L12").Value = -1153933
$ws.Range("M12").Value = -1226116
# Row 13
$ws.Range("D13").Value = 587812
$ws.Range("E13").Value = 306402
$ws.Range("F13").Value = 1111886
$ws.Range("G13").Value = 1881881
$ws.Range("H13").Value = 894531
$ws.Range("I13").Value = 362657
$ws.Range("J13").Value = 1732455
$ws.Range("K13").Value = 1594945
$ws.Range("L13").Value = 1217118
$ws.Range("M13").Value = 1045479
# Row 14
$ws.Range("D14").Value = -23409
$ws.Range("E14").Value = -23596
$ws.Range("F14").Value = -72740
$ws.Range("G14").Value = -82626
$ws.Range("H14").Value = -62266
$ws.Range("I14").Value = -67411
$ws.Range("J14").Value = -121283
$ws.Range("K14").Value = -79427
$ws.Range("L14").Value = -84957
$ws.Range("M14").Value = -209952
# Row 16
$ws.Range("D16").Value = 3065
$ws.Range("E16").Value = -7120
$ws.Range("F16").Value = 2236
$ws.Range("G16").Value = 12436
$ws.Range("H16").Value = 1957
$ws.Range("I16").Value = -12891
$ws.Range("J16").Value = 1284
$ws.Range("K16").Value = 55888
$ws.Range("L16").Value = 5392
$ws.Range("M16").Value = 305847
# Row 17
$ws.Range("D17").Value = 567468
$ws.Range("E17").Value = 275686
$ws.Range("F17").Value = 1041382
$ws.Range("G17").Value = 1811691
$ws.Range("H17").Value = 834222
$ws.Range("I17").Value = 282355
$ws.Range("J17").Value = 1612456
$ws.Range("K17").Value = 1571406
$ws.Range("L17").Value = 1137553
$ws.Range("M17").Value = 1141374
# Row 18
$ws.Range("D18").Value = -1083
$ws.Range("E18").Value = 51
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = -4919
$ws.Range("L18").Value = -7500
$ws.Range("M18").Value = -17713
# Row 19
$ws.Range("D19").Value = 64893
$ws.Range("E19").Value = 406034
$ws.Range("F19").Value = 641286
$ws.Range("G19").Value = 64522
$ws.Range("H19").Value = 113809
$ws.Range("I19").Value = 149991
$ws.Range("J19").Value = 116913
$ws.Range("K19").Value = 136596
$ws.Range("L19").Value = 104814
$ws.Range("M19").Value = 240532
# Row 20
$ws.Range("D20").Value = 631278
$ws.Range("E20").Value = 681771
$ws.Range("F20").Value = 1682668
$ws.Range("G20").Value = 1876213
$ws.Range("H20").Value = 948031
$ws.Range("I20").Value = 432346
$ws.Range("J20").Value = 1729369
$ws.Range("K20").Value = 1703083
$ws.Range("L20").Value = 1234867
$ws.Range("M20").Value = 1364193
# Row 21
$ws.Range("D21").Value = -77256
$ws.Range("E21").Value = 149239
$ws.Range("F21").Value = -213916
$ws.Range("G21").Value = -364764
$ws.Range("H21").Value = -93957
$ws.Range("I21").Value = 232709
$ws.Range("J21").Value = -250593
$ws.Range("K21").Value = -258558
$ws.Range("L21").Value = -149853
$ws.Range("M21").Value = 170906
# Row 22
$ws.Range("D22").Value = 554022
$ws.Range("E22").Value = 831010
$ws.Range("F22").Value = 1468752
$ws.Range("G22").Value = 1511449
$ws.Range("H22").Value = 854074
$ws.Range("I22").Value = 665055
$ws.Range("J22").Value = 1478776
$ws.Range("K22").Value = 1444525
$ws.Range("L22").Value = 1085014
$ws.Range("M22").Value = 1535099
# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 6448
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 5543
# Row 24
$ws.Range("D24").Value = 554022
$ws.Range("E24").Value = 831010
$ws.Range("F24").Value = 1468752
$ws.Range("G24").Value = 1511449
$ws.Range("H24").Value = 854074
$ws.Range("I24").Value = 665055
$ws.Range("J24").Value = 1478776
$ws.Range("K24").Value = 1444525
$ws.Range("L24").Value = 1085014
$ws.Range("M24").Value = 1535099
# Row 25
$ws.Range("D25").Value = 791
$ws.Range("E25").Value = 1187
$ws.Range("F25").Value = 2098
$ws.Range("G25").Value = 2159
$ws.Range("H25").Value = 1220
$ws.Range("I25").Value = 950
$ws.Range("J25").Value = 2113
$ws.Range("K25").Value = 2064
$ws.Range("L25").Value = 1085
$ws.Range("M25").Value = 1535
# Row 26
$ws.Range("D26").Value = 700000
$ws.Range("E26").Value = 700000
$ws.Range("F26").Value = 700000
$ws.Range("G26").Value = 700000
$ws.Range("H26").Value = 700000
$ws.Range("I26").Value = 700000
$ws.Range("J26").Value = 700000
$ws.Range("K26").Value = 700000
$ws.Range("L26").Value = 1000000
$ws.Range("M26").Value = 1000000
# Row 27
$ws.Range("D27").Value = 554
$ws.Range("E27").Value = 831
$ws.Range("F27").Value = 1469
$ws.Range("G27").Value = 1511
$ws.Range("H27").Value = 854
$ws.Range("I27").Value = 665
$ws.Range("J27").Value = 1479
$ws.Range("K27").Value = 1445
$ws.Range("L27").Value = 1085
$ws.Range("M27").Value = 1535
